$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "掌阅科技"
$ws.Range("B2").Value = "华胜天成"
$ws.Range("C2").Value = "华胜天成"

$ws.Range("A3").Value = "光线传媒"
$ws.Range("B3").Value = "深科技"
$ws.Range("C3").Value = "巨力索具"

$ws.Range("A4").Value = "大位科技"
$ws.Range("B4").Value = "光线传媒"
$ws.Range("C4").Value = "汉缆股份"

$ws.Range("A5").Value = "华胜天成"
$ws.Range("B5").Value = "利欧股份"
$ws.Range("C5").Value = "嘉美包装"

$ws.Range("A6").Value = "巨力索具"
$ws.Range("B6").Value = "掌阅科技"
$ws.Range("C6").Value = "利欧股份"

$ws.Range("A7").Value = "利欧股份"
$ws.Range("B7").Value = "天奇股份"
$ws.Range("C7").Value = "掌阅科技"

$ws.Range("A8").Value = "深科技"
$ws.Range("B8").Value = "巨力索具"
$ws.Range("C8").Value = "博纳影业"

$ws.Range("A9").Value = "汉缆股份"
$ws.Range("B9").Value = "大位科技"
$ws.Range("C9").Value = "光线传媒"

$ws.Range("A10").Value = "天奇股份"
$ws.Range("B10").Value = "兴民智通"
$ws.Range("C10").Value = "协鑫集成"

$ws.Range("A11").Value = "嘉美包装"
$ws.Range("B11").Value = "汉缆股份"
$ws.Range("C11").Value = "大位科技"

$ws.Range("A12").Value = "捷成股份"
$ws.Range("B12").Value = "东方财富"
$ws.Range("C12").Value = "天奇股份"

$ws.Range("A13").Value = "协鑫集成"
$ws.Range("B13").Value = "贵州茅台"
$ws.Range("C13").Value = "紫金矿业"

$ws.Range("A14").Value = "浙江世宝"
$ws.Range("B14").Value = "航发动力"
$ws.Range("C14").Value = "深科技"

$ws.Range("A15").Value = "博纳影业"
$ws.Range("B15").Value = "嘉美包装"
$ws.Range("C15").Value = "万向钱潮"

$ws.Range("A16").Value = "航发动力"
$ws.Range("B16").Value = "协鑫集成"
$ws.Range("C16").Value = "风语筑"

$ws.Range("A17").Value = "国安股份"
$ws.Range("B17").Value = "亚星锚链"
$ws.Range("C17").Value = "航天发展"

$ws.Range("A18").Value = "兴民智通"
$ws.Range("B18").Value = "博纳影业"
$ws.Range("C18").Value = "国安股份"

$ws.Range("A19").Value = "风语筑"
$ws.Range("B19").Value = "天汽模"
$ws.Range("C19").Value = "二六三"

$ws.Range("A20").Value = "东方财富"
$ws.Range("B20").Value = "国安股份"
$ws.Range("C20").Value = "五洲新春"

$ws.Range("A21").Value = "贵州茅台"
$ws.Range("B21").Value = "风语筑"
$ws.Range("C21").Value = "易天股份"
